# Generate Report for Handback
# Updates the handback-status workbook with the freshly generated
# handback file names / UUIDs and timestamps for this CI run.

$wb = $excel.ActiveWorkbook

$newId1 = "d4bc9dbf-0870-47a8-963f-5aede2dd074f"
$newId2 = "ffffa1516d6c-1b7a-44e3-a722-b57c9951bffb"

$newXlfBase = "d4bc9dbf-0870-47a8-963f-5aede2dd074f.3d9250b61cab37357bbae60693ea6debfc0d761b"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G2").Value = "2016-08-22 17:05:25"
$wsOverview.Range("G3").Value = "2016-08-22 17:05:25"

# Hyperlinks in document order: B2, B3.
# NOTE: must iterate with foreach (indexing via .Item(n) + mutate duplicates
# the hyperlink entry in this host), updating TextToDisplay in place so the
# r:id relationship is preserved.
$overviewDisplays = @("e2e\$newId1.md", "e2e\$newId2.md")
$i = 0
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $overviewDisplays[$i]
    $i = $i + 1
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("I3").Value = "$newId2.md"

$wsZhCn.Range("G2").Value = "$newXlfBase.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$newXlfBase.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "$newXlfBase.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "$newXlfBase.zh-cn.xlf"

$wsZhCn.Range("H2").Value = "2016-08-22 17:05:20"
$wsZhCn.Range("H3").Value = "2016-08-22 17:05:20"
$wsZhCn.Range("K2").Value = "2016-08-22 17:05:37"
$wsZhCn.Range("K3").Value = "2016-08-22 17:05:37"

# Hyperlinks in document order: A2, I2, A3, I3
$zhCnDisplays = @("$newId1.md", "$newId1.md", "$newId2.md", "$newId2.md")
$i = 0
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $zhCnDisplays[$i]
    $i = $i + 1
}

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("I3").Value = "$newId2.md"

$wsDeDe.Range("G2").Value = "$newXlfBase.de-de.xlf"
$wsDeDe.Range("J2").Value = "$newXlfBase.de-de.xlf"
$wsDeDe.Range("G3").Value = "$newXlfBase.de-de.xlf"
$wsDeDe.Range("J3").Value = "$newXlfBase.de-de.xlf"

$wsDeDe.Range("H2").Value = "2016-08-22 17:05:25"
$wsDeDe.Range("H3").Value = "2016-08-22 17:05:25"
$wsDeDe.Range("K2").Value = "2016-08-22 17:05:44"
$wsDeDe.Range("K3").Value = "2016-08-22 17:05:44"

# Hyperlinks in document order: A2, I2, A3, I3
$deDeDisplays = @("$newId1.md", "$newId1.md", "$newId2.md", "$newId2.md")
$i = 0
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $deDeDisplays[$i]
    $i = $i + 1
}
